$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.689.91'
$ws.Range("E2").Value = '  +5.69%  '
$ws.Range("D3").Value = '3.466.32'
$ws.Range("E3").Value = '  +3.29%  '
$ws.Range("E4").Value = '  +0.26%  '
$rng = $ws.Range("D5")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '410.66'
$rng.Style = $origStyle
$ws.Range("E5").Value = '  -1.02%  '
$rng = $ws.Range("D6")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '128.86'
$rng.Style = $origStyle
$ws.Range("E6").Value = '  +15.19%  '
$ws.Range("D7").Value = '3.460.97'
$ws.Range("E7").Value = '  +3.27%  '
$rng = $ws.Range("D8")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '0.594'
$rng.Style = $origStyle
$ws.Range("E8").Value = '  +1.19%  '
$ws.Range("E9").Value = '  +0.05%  '
$ws.Range("E10").Value = '  +8.72%  '
$ws.Range("E11").Value = '  +30.70%  '
$rng = $ws.Range("D12")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '43.67'
$rng.Style = $origStyle
$ws.Range("E12").Value = '  +8.61%  '
$ws.Range("E13").Value = '  +0.17%  '
$ws.Range("D14").Value = '4.008.73'
$ws.Range("E14").Value = '  +3.68%  '
$ws.Range("E15").Value = '  +2.80%  '
$rng = $ws.Range("D16")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '20.22'
$rng.Style = $origStyle
$ws.Range("E16").Value = '  +3.52%  '
$ws.Range("D17").Value = '3.503.08'
$ws.Range("E17").Value = '  +3.64%  '
$ws.Range("D18").Value = '62.550.53'
$ws.Range("E18").Value = '  +5.91%  '
$ws.Range("E19").Value = '  +0.69%  '
$rng = $ws.Range("D20")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '10.95'
$rng.Style = $origStyle
$ws.Range("E20").Value = '  -0.20%  '
$rng = $ws.Range("D21")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '0.0000137'
$rng.Style = $origStyle
$ws.Range("E21").Value = '  +25.25%  '
$rng = $ws.Range("D22")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '3.37'
$rng.Style = $origStyle
$ws.Range("E22").Value = '  +0.00%  '
$rng = $ws.Range("D23")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '13.22'
$rng.Style = $origStyle
$ws.Range("E23").Value = '  +1.05%  '
$rng = $ws.Range("D24")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '82.09'
$rng.Style = $origStyle
$ws.Range("E24").Value = '  +8.98%  '
$rng = $ws.Range("D25")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '312.58'
$rng.Style = $origStyle
$ws.Range("E25").Value = '  +2.84%  '
$rng = $ws.Range("D26")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '3.19'
$rng.Style = $origStyle
$ws.Range("E26").Value = '  -1.60%  '
$rng = $ws.Range("D27")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '30.35'
$rng.Style = $origStyle
$ws.Range("E27").Value = '  +5.69%  '
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$rng = $ws.Range("D28")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '7.85'
$rng.Style = $origStyle
$ws.Range("E28").Value = '  +6.00%  '
$ws.Range("B29").Value = 'Filecoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$rng = $ws.Range("D29")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '8.09'
$rng.Style = $origStyle
$ws.Range("E29").Value = '  +0.94%  '
$rng = $ws.Range("D30")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '0.121'
$rng.Style = $origStyle
$ws.Range("E30").Value = '  +7.35%  '
$ws.Range("E31").Value = '  +3.77%  '
$rng = $ws.Range("D32")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '4.36'
$rng.Style = $origStyle
$ws.Range("E32").Value = '  -2.30%  '
$rng = $ws.Range("D33")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '44.97'
$rng.Style = $origStyle
$ws.Range("E33").Value = '  +12.02%  '
$rng = $ws.Range("D34")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '2.70'
$rng.Style = $origStyle
$ws.Range("E34").Value = '  +26.49%  '
$ws.Range("E35").Value = '  +4.30%  '
$ws.Range("E36").Value = '  +0.10%  '
$rng = $ws.Range("D37")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '0.0492'
$rng.Style = $origStyle
$ws.Range("E37").Value = '  -8.79%  '
$rng = $ws.Range("D38")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '52.73'
$rng.Style = $origStyle
$ws.Range("E38").Value = '  +0.94%  '
$ws.Range("E39").Value = '  +2.14%  '
$ws.Range("E40").Value = '  -0.15%  '
$ws.Range("E41").Value = '  -6.69%  '
$ws.Range("E42").Value = '  +3.80%  '
$rng = $ws.Range("D43")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '18.14'
$rng.Style = $origStyle
$ws.Range("E43").Value = '  +7.24%  '
$rng = $ws.Range("D44")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '137.72'
$rng.Style = $origStyle
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("E45").Value = '  +2.46%  '
$rng = $ws.Range("D46")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '0.291'
$rng.Style = $origStyle
$ws.Range("E46").Value = '  +3.90%  '
$rng = $ws.Range("D47")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '3.99'
$rng.Style = $origStyle
$ws.Range("E47").Value = '  -0.37%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$rng = $ws.Range("D48")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '22.64'
$rng.Style = $origStyle
$ws.Range("E48").Value = '  -0.69%  '
$ws.Range("B49").Value = 'WEMIXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$rng = $ws.Range("D49")
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '2.25'
$rng.Style = $origStyle
$ws.Range("E49").Value = '  -1.23%  '
$ws.Range("D50").Value = '2.244.61'
$ws.Range("E50").Value = '  +1.76%  '
$ws.Range("D51").Value = '3.807.17'
$ws.Range("E51").Value = '  +3.83%  '
